# This script applies the edit described by the diff:
#  - Removes w:proofErr spell/gram-check markers and merges the runs
#    that had been split around them back into single runs.
#  - Moves the _GoBack bookmark from its own (empty) paragraph right
#    after the Fig. 3 picture to the end of the last bullet item in the
#    "Goals of this semester" list, leaving an empty paragraph behind.
#
# Each affected paragraph is replaced wholesale via Range.InsertXML on
# just that paragraph's Range, which performs a precise, local XML
# substitution without perturbing the rest of the document.

$d = $word.ActiveDocument

# "Iterate data sources..." bullet - gains the _GoBack bookmark at its end
$p = $d.Paragraphs.Item(35)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" w:rsidR="00D92D0D" w:rsidRPr="001804BE" w:rsidRDefault="00D92D0D" w:rsidP="00D92D0D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Iterate data sources to collect data. For example, evaluate software engineering documents to collect relevant information. (5)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@)

# "Model Training (...)" bullet - merge split run around "size"
$p = $d.Paragraphs.Item(30)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" w:rsidR="0099337C" w:rsidRDefault="0099337C" w:rsidP="0099337C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Model Training (this part may require domain knowledge of software engineering and software size estimation).</w:t></w:r></w:p>
'@)

# "Improve the APIs..." bullet - merge split run around "apis"
$p = $d.Paragraphs.Item(23)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" w:rsidR="0099337C" w:rsidRDefault="0099337C" w:rsidP="0099337C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Improve the APIs, for example, check validity of the uploaded files and parameters, restructure the data base query functions for better efficiency, put comments on the apis for its readability. (3)</w:t></w:r></w:p>
'@)

# Paragraph 18: paragraph that only held the _GoBack bookmark becomes empty
$p = $d.Paragraphs.Item(18)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" w:rsidR="0099337C" w:rsidRDefault="0099337C" w:rsidP="00747D4D"/>
'@)

# Paragraph 9: "Fig. 1. UMLx system architecture." caption - merge split run
$p = $d.Paragraphs.Item(9)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" w:rsidR="00747D4D" w:rsidRDefault="007856A1" w:rsidP="00747D4D"><w:r><w:t>Fig. 1. UMLx system architecture.</w:t></w:r></w:p>
'@)

# Paragraph 7: merge the two long runs that had spell/gram-check splits
$p = $d.Paragraphs.Item(7)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" w:rsidR="00902D46" w:rsidRDefault="00747D4D" w:rsidP="00747D4D"><w:r><w:t>For the purpose of standardizing our proposed approach to analysing UML diagrams and improving its efficiency, the proposed procedure and algorithms are implemented as a tool called UMLx, which is an abbreviation for UML analytics. This web-based UML diagrams analytic tool is implemented with Nodejs/MongoDB for the back-end logic and data storage, HTML/Javascript/CSS for the front end, and R scripts for statistical analyses.</w:t></w:r><w:r w:rsidR="007856A1"><w:t xml:space="preserve"> The</w:t></w:r><w:r><w:t xml:space="preserve"> Model File Parser module is currently implemented based on the format of the XML files exported from Enterprise Architect, to extend the use of the tool to other UML modeling tools, for example, Visual paradigm, an individual parser needs to be implemented and integrated to be able to analyse specific types of XML files.</w:t></w:r></w:p>
'@)

# Paragraph 4: "The structure of UMLx..." - merge split run + remove gramStart/End around "Matching"
$p = $d.Paragraphs.Item(4)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" w:rsidR="00747D4D" w:rsidRDefault="00747D4D" w:rsidP="00747D4D"><w:r><w:t xml:space="preserve">The structure of UMLx is shown in Fig. </w:t></w:r><w:r w:rsidR="007856A1"><w:t>1</w:t></w:r><w:r><w:t xml:space="preserve">, which comprises 7 modules to support the analyses. Specifically, Model File Parser module is responsible for parsing the input XML files into the directed graphs and associating nodes with the corresponding UML diagram elements. Element Consolidation module is to consolidate the elements that share the same name across different diagrams. Path Profiler module traverses the directed graphs to identify the paths that are the potential transactions. Pattern Tree Builder module creates a tree structure for the input patterns, which will be used later in pattern matching process. Pattern </w:t></w:r><w:r w:rsidR="002C5A85"><w:t>Matching</w:t></w:r><w:r><w:t xml:space="preserve"> module matches each of identified paths against the pattern tree and associates operational characteristics with the transaction when a path is matched and determined as a transaction. Identified transactions are input to Performance Evaluator module to be compared with the empirically understood transactions in order to calculate identification and categorization accuracy. Also the identified transactions are input to Transaction Analyzer module to perform calculations according to the defined metrics. The end users are able to review the results at different stages of the procedure.</w:t></w:r><w:r w:rsidR="00EB251A"><w:t xml:space="preserve"> A layered architecture view also provided in Fig. 2.</w:t></w:r></w:p>
'@)

# Paragraph 1: title - merge split UMLx run
$p = $d.Paragraphs.Item(1)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" w:rsidR="00A82D77" w:rsidRDefault="00A82D77" w:rsidP="00A82D77"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Current architecture and implementation for UMLx</w:t></w:r></w:p>
'@)
